$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts "1.003" etc. into a
# number. Style is reset back to Normal afterwards so no stray style diff
# is introduced.
$ws.Range("D2").Value = '26.950.30'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '1.817.58'
$ws.Range("E3").Value = '  +0.50%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4692'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.85%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3664'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07351'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8723'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.26'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.89%  '
$ws.Range("D12").Value = '1.817.03'
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.404'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07115'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.509'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.06%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.35'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008702'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("D21").Value = '26.967.41'
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.293'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.41%  '
$ws.Range("E23").Value = '  +0.74%  '
$ws.Range("D24").Value = '2.046.13'
$ws.Range("E24").Value = '  +0.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.895'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.85'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.50%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.00%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.148'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.28%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.252'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.04'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08906'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.10%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7589'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("E33").Value = '  +0.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.498'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.910'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.35%  '
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("E37").Value = '  -0.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05289'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.69%  '
$ws.Range("E39").Value = '  -0.87%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.970'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.382'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.88%  '
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5285'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.20%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '7.152'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.429'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.07%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4868'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.29%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.47'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.79%  '
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '103.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.663'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06297'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.25%  '
